{"js": "// Insert a new \"List Bullet\" paragraph with the two instructors, right\n// after the \"Docente(s) Respons\u00e1vel(eis)\" heading paragraph and before\n// the \"Programa resumido\" heading paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Docente(s) Respons\u00e1vel(eis)\") !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the 'Docente(s) Respons\u00e1vel(eis)' paragraph\");\n}\n\n// A zero-length range right after the heading paragraph's end-of-paragraph\n// mark -- insertOoxml here adds a brand-new paragraph between it and the\n// paragraph that follows (\"Programa resumido\").\nconst afterRange = target.getRange(Word.RangeLocation.after);\n\nconst ooxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n          '<w:body>' +\n            '<w:p>' +\n              '<w:pPr><w:pStyle w:val=\"ListBullet\"/></w:pPr>' +\n              '<w:r><w:t>5111420 - Talita Martins Lacerda</w:t><w:br/></w:r>' +\n              '<w:r><w:t>8853480 - Tatiane da Franca Silva</w:t></w:r>' +\n            '</w:p>' +\n          '</w:body>' +\n        '</w:document>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>';\n\nafterRange.insertOoxml(ooxml);\nawait context.sync();\n", "ps1": "# Insert a new \"List Bullet\" paragraph with the two instructors, right\n# after the \"Docente(s) Respons\u00e1vel(eis)\" heading paragraph and before\n# the \"Programa resumido\" heading paragraph.\n\n$d = $word.ActiveDocument\n\n# Locate the \"Docente(s) Respons\u00e1vel(eis)\" paragraph by its text (plain\n# substring check -- avoids any regex-metacharacter escaping pitfalls).\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Contains(\"Docente(s) Respons\u00e1vel(eis)\")) {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find the 'Docente(s) Respons\u00e1vel(eis)' paragraph\"\n}\n\n# Position a zero-length range right before the paragraph mark that ends\n# the \"Docente(s) Respons\u00e1vel(eis)\" paragraph, so the new OOXML is\n# inserted as a brand-new paragraph between it and the next one.\n$insertAt = $target.Range.End - 1\n$insertionRange = $d.Range($insertAt, $insertAt)\n\n$newParagraphXml = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"ListBullet\"/>\n            </w:pPr>\n            <w:r>\n              <w:t>5111420 - Talita Martins Lacerda</w:t>\n              <w:br/>\n            </w:r>\n            <w:r>\n              <w:t>8853480 - Tatiane da Franca Silva</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n'@\n\n$insertionRange.InsertXML($newParagraphXml)\n"}
